$d = $word.ActiveDocument

# The title paragraph reads "Descrição dos Casos de Uso". The bold,
# 48-half-point (24pt) run ending in "Casos" sits immediately before
# the hidden "_GoBack" bookmark. The edit adds two plain spaces
# ("  "), in the same bold/24pt formatting, right after "Casos" and
# right before that bookmark.

$rng = $null

# Prefer anchoring on the "_GoBack" bookmark - it marks the exact
# insertion point used by the edit and is more precise than a text
# search.
if ($d.Bookmarks.Exists("_GoBack")) {
    $rng = $d.Bookmarks.Item("_GoBack").Range
    $rng.Collapse(1)   # wdCollapseStart
}
else {
    # Fallback: locate the end of "Casos" via Find.
    $rng = $d.Content
    $found = $rng.Find.Execute("Casos", $true, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Collapse(0)   # wdCollapseEnd
    }
    else {
        $rng = $null
    }
}

if ($rng -ne $null) {
    # Match the formatting already used by the surrounding title
    # text: bold, 24pt (sz/szCs = 48 half-points). The inserted text
    # naturally inherits the Portuguese (Brazil) language already in
    # effect at this point in the paragraph.
    $rng.Font.Bold = $true
    $rng.Font.Size = 24

    $rng.InsertAfter("  ")
}
